$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 863.4
$ws.Range("I2").Value = 105.666664
$ws.Range("K2").Value = 105.666664
$ws.Range("M2").Value = 7.333336000000003

$ws.Range("H41").Value = 204.85715
$ws.Range("I41").Value = 424.5
$ws.Range("J41").Value = 117
$ws.Range("K41").Value = 424.5
$ws.Range("L41").Value = 117
$ws.Range("M41").Value = 15.5
$ws.Range("N41").Value = -997

$ws.Range("H68").Value = 69420
$ws.Range("J68").Value = 69420
$ws.Range("L68").Value = 69420
$ws.Range("N68").Value = -70918

$ws.Range("H71").Value = 69420
$ws.Range("J71").Value = 69420
$ws.Range("L71").Value = 208260
$ws.Range("N71").Value = -215748

$ws.Range("H98").Value = 10706.394
$ws.Range("I98").Value = 1583.3125
$ws.Range("K98").Value = 1583.3125
$ws.Range("M98").Value = -85.3125

$ws.Range("H105").Value = 26835
$ws.Range("J105").Value = 26835
$ws.Range("L105").Value = 26835
$ws.Range("N105").Value = -33823

$ws.Range("H122").Value = 10706.394
$ws.Range("I122").Value = 1583.3125
$ws.Range("K122").Value = 4749.9375
$ws.Range("M122").Value = -2299.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1422.9
$ws.Range("I2").Value = 1072.6428
$ws.Range("J2").Value = 2240.1667
$ws.Range("K2").Value = 1072.6428
$ws.Range("L2").Value = 2240.1667
$ws.Range("M2").Value = -959.6428000000001
$ws.Range("N2").Value = -2466.1667

$ws.Range("H45").Value = 3725.182
$ws.Range("I45").Value = 2082.2
$ws.Range("K45").Value = 2082.2
$ws.Range("M45").Value = -1705.2

$ws.Range("H102").Value = 6446.0625
$ws.Range("I102").Value = 6446.0625
$ws.Range("K102").Value = 6446.0625
$ws.Range("M102").Value = -4824.0625

$ws.Range("H110").Value = 834.5333000000001
$ws.Range("I110").Value = 869.1667
$ws.Range("J110").Value = 696
$ws.Range("K110").Value = 869.1667
$ws.Range("L110").Value = 696
$ws.Range("M110").Value = 1175.8333
$ws.Range("N110").Value = -4786

$ws.Range("H116").Value = 1422.9
$ws.Range("I116").Value = 1072.6428
$ws.Range("J116").Value = 2240.1667
$ws.Range("K116").Value = 1072.6428
$ws.Range("L116").Value = 2240.1667
$ws.Range("M116").Value = 1221.3572
$ws.Range("N116").Value = -6828.1667

$ws.Range("H132").Value = 3951
$ws.Range("I132").Value = 3853.225
$ws.Range("K132").Value = 11559.675
$ws.Range("M132").Value = -9029.674999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1422.9
$ws.Range("I3").Value = 1072.6428
$ws.Range("J3").Value = 2240.1667
$ws.Range("K3").Value = 1072.6428
$ws.Range("L3").Value = 2240.1667
$ws.Range("M3").Value = -958.6428000000001
$ws.Range("N3").Value = -2468.1667

$ws.Range("H36").Value = 4975
$ws.Range("J36").Value = 13020
$ws.Range("L36").Value = 13020
$ws.Range("N36").Value = -14088

$ws.Range("H64").Value = 1444.5
$ws.Range("I64").Value = 1225.6666
$ws.Range("K64").Value = 1225.6666
$ws.Range("M64").Value = -1000.6666

$ws.Range("H67").Value = 1444.5
$ws.Range("I67").Value = 1225.6666
$ws.Range("K67").Value = 1225.6666
$ws.Range("M67").Value = -445.6666

$ws.Range("H76").Value = 314
$ws.Range("J76").Value = 314
$ws.Range("L76").Value = 314
$ws.Range("N76").Value = -944

$ws.Range("H79").Value = 314
$ws.Range("J79").Value = 314
$ws.Range("L79").Value = 314
$ws.Range("N79").Value = -2498

$ws.Range("H94").Value = 477.6154
$ws.Range("I94").Value = 323.78946
$ws.Range("J94").Value = 895.1429000000001
$ws.Range("K94").Value = 323.78946
$ws.Range("L94").Value = 895.1429000000001
$ws.Range("M94").Value = 127.21054
$ws.Range("N94").Value = -1797.1429

$ws.Range("H99").Value = 1110.875
$ws.Range("I99").Value = 1110.875
$ws.Range("K99").Value = 1110.875
$ws.Range("M99").Value = 387.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 302.1
$ws.Range("I7").Value = 373
$ws.Range("K7").Value = 373
$ws.Range("M7").Value = -260

$ws.Range("H16").Value = 495.8
$ws.Range("I16").Value = 582.25
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = 582.25
$ws.Range("L16").Value = 150
$ws.Range("M16").Value = -295.25
$ws.Range("N16").Value = -724

$ws.Range("H105").Value = 1751.7273
$ws.Range("I105").Value = 1807.7778
$ws.Range("K105").Value = 1807.7778
$ws.Range("M105").Value = -60.77780000000007

$ws.Range("H107").Value = 324.83334
$ws.Range("I107").Value = 149.8
$ws.Range("K107").Value = 149.8
$ws.Range("M107").Value = 1770.2

$ws.Range("H113").Value = 495.8
$ws.Range("I113").Value = 582.25
$ws.Range("J113").Value = 150
$ws.Range("K113").Value = 582.25
$ws.Range("L113").Value = 150
$ws.Range("M113").Value = 1587.75
$ws.Range("N113").Value = -4490

$ws.Range("H132").Value = 2687
$ws.Range("I132").Value = 2296.8965
$ws.Range("K132").Value = 6890.689499999999
$ws.Range("M132").Value = -4360.689499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 604.9375
$ws.Range("I5").Value = 578.6667
$ws.Range("J5").Value = 999
$ws.Range("K5").Value = 1736.0001
$ws.Range("L5").Value = 2997
$ws.Range("M5").Value = -1624.0001
$ws.Range("N5").Value = -3221

$ws.Range("H12").Value = 86.818184
$ws.Range("J12").Value = 66.71429000000001
$ws.Range("L12").Value = 200.14287
$ws.Range("N12").Value = -546.14287

$ws.Range("H121").Value = 79601.53
$ws.Range("J121").Value = 104291
$ws.Range("L121").Value = 312873
$ws.Range("N121").Value = -315493

$ws.Range("H131").Value = 1471.4736
$ws.Range("I131").Value = 952.5454999999999
$ws.Range("J131").Value = 2185
$ws.Range("K131").Value = 2857.6365
$ws.Range("L131").Value = 6555
$ws.Range("M131").Value = 2182.3635
$ws.Range("N131").Value = -16635

$ws.Range("H135").Value = 604.9375
$ws.Range("I135").Value = 578.6667
$ws.Range("J135").Value = 999
$ws.Range("K135").Value = 5208.0003
$ws.Range("L135").Value = 8991
$ws.Range("M135").Value = -2673.0003
$ws.Range("N135").Value = -14061

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9392.4
$ws.Range("I132").Value = 4674.5
$ws.Range("J132").Value = 10571.875
$ws.Range("K132").Value = 14023.5
$ws.Range("L132").Value = 31715.625
$ws.Range("M132").Value = -11493.5
$ws.Range("N132").Value = -36775.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 95000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 95000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 95000
$ws.Range("N87").Value = -97246
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 95000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 95000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 285000
$ws.Range("N90").Value = -296232
$ws.Range("M90").ClearContents()

$ws.Range("H93").Value = 1386.1786
$ws.Range("I93").Value = 1266.1
$ws.Range("K93").Value = 1266.1
$ws.Range("M93").Value = -18.09999999999991

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 49228.625
$ws.Range("I70").Value = 24500
$ws.Range("J70").Value = 57471.5
$ws.Range("K70").Value = 24500
$ws.Range("L70").Value = 57471.5
$ws.Range("M70").Value = -24185
$ws.Range("N70").Value = -58101.5

$ws.Range("H73").Value = 49228.625
$ws.Range("I73").Value = 24500
$ws.Range("J73").Value = 57471.5
$ws.Range("K73").Value = 24500
$ws.Range("L73").Value = 57471.5
$ws.Range("M73").Value = -23408
$ws.Range("N73").Value = -59655.5

$ws.Range("H132").Value = 2044.0385
$ws.Range("I132").Value = 1836.5454
$ws.Range("J132").Value = 3185.25
$ws.Range("K132").Value = 5509.6362
$ws.Range("L132").Value = 9555.75
$ws.Range("M132").Value = -2979.6362
$ws.Range("N132").Value = -14615.75
